# Update the "dSF" column (F) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = 1
    5  = 1
    6  = 3
    7  = 8
    8  = -4
    9  = -3
    10 = 2
    11 = 3
    12 = -3
    14 = 3
    15 = 5
    16 = 6
    17 = -2
    18 = 7
    19 = -2
    20 = 2
    21 = -4
    22 = -2
    23 = 1
    25 = -7
    26 = 4
    28 = 2
    29 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
